# generate_ruml.py dependency diagram: shift the "flow arrow" marker cells
# one column to the right for the CarPollutionPermit / Vehicle / CheckFraud
# branches of the tree (and touch up rows 5, 12, 15, 18, 21 so the arrows
# land in their new, correct columns).
#
# These marker cells share one cell style (a bordered cell used to draw
# the tree lines, taken here from A6 - a cell that is never itself moved
# or cleared by this script, so it stays a valid format source for the
# whole run). Moving a marker to a column that doesn't have that border
# yet means copying the format in; the column that a marker vacates for
# good must be fully cleared (value *and* format) so no stray
# bordered-but-empty cell is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fmtSource = $ws.Range("A6")

function Copy-CellFormat {
    param($targetAddress)
    $fmtSource.Copy() | Out-Null
    $ws.Range($targetAddress).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

# --- Rows 2-4: move the marker from column A to column C ---------------
Copy-CellFormat "C2"
$ws.Range("C2").Value = $ws.Range("A2").Value()
$ws.Range("A2").Clear()

Copy-CellFormat "C3"
$ws.Range("A3").Clear()

Copy-CellFormat "C4"
$ws.Range("A4").Clear()

# --- Row 5: marker moves from B to A; C gains the (empty) bordered cell
Copy-CellFormat "C5"
$ws.Range("A5").Value = $ws.Range("B5").Value()
$ws.Range("B5").Clear()

# --- Rows 6-11: move the marker from column B to column C --------------
Copy-CellFormat "C6"
$ws.Range("B6").Clear()

Copy-CellFormat "C7"
$ws.Range("B7").Clear()

Copy-CellFormat "C8"
$ws.Range("B8").Clear()

Copy-CellFormat "C9"
$ws.Range("B9").Clear()

Copy-CellFormat "C10"
$ws.Range("B10").Clear()

Copy-CellFormat "C11"
$ws.Range("B11").Clear()

# --- Row 12: marker moves from D to B; C gains the (empty) bordered cell
Copy-CellFormat "C12"
$ws.Range("B12").Value = $ws.Range("D12").Value()
$ws.Range("D12").Clear()

# --- Rows 13-14, 16-17: move the marker from column D to column C ------
Copy-CellFormat "C13"
$ws.Range("D13").Clear()

Copy-CellFormat "C14"
$ws.Range("D14").Clear()

Copy-CellFormat "C16"
$ws.Range("D16").Clear()

Copy-CellFormat "C17"
$ws.Range("D17").Clear()

# --- Row 15: A's marker (←) value clears but the bordered cell stays in
#             place; B gets D's old marker (◁); C (newly-formatted) takes
#             A's old marker (←); D is fully vacated -------------------
Copy-CellFormat "C15"
$oldA15 = $ws.Range("A15").Value()
$oldD15 = $ws.Range("D15").Value()
$ws.Range("C15").Value = $oldA15
$ws.Range("B15").Value = $oldD15
$ws.Range("A15").ClearContents()
$ws.Range("D15").Clear()

# --- Row 18: marker moves from column C to column D (both cells already
#             have the bordered style, so only the value needs to move) --
$ws.Range("D18").Value = $ws.Range("C18").Value()
$ws.Range("C18").ClearContents()

# --- Row 21: A and D markers swap their arrow glyphs (B and C unchanged,
#             all four cells keep their bordered style) -----------------
$ws.Range("A21").Value = "◁"
$ws.Range("D21").Value = "←"

Write-Host "generate_ruml.py dependency diagram updated"
